{"js": "// The template document contains a paragraph whose text reads \"{m: null}\".\n// That paragraph has two runs: \"{m: \" and \"null}\". The edit splits the\n// second run (\"null}\") into two runs: \"null\" and \"}\" (the latter flagged\n// with xml:space=\"preserve\"), without touching anything else.\n\nconst body = context.document.body;\n\n// Locate the run of text \"null}\" (there is exactly one occurrence in the\n// document, inside the \"{m: null}\" paragraph).\nconst results = body.search(\"null}\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find text \"null}\" to split.');\n}\n\nconst target = results.items[0];\n\n// Replace the matched range with explicit OOXML that keeps the original\n// run's formatting/identity for \"null\" (preserving its w:rsidR attribute)\n// and introduces a brand-new run for the trailing \"}\" character, exactly\n// mirroring a manual edit/retype of the closing brace in Word.\nconst replacementOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p>' +\n  '<w:r w:rsidR=\"00D500A1\"><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>null</w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  '</w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntarget.insertOoxml(replacementOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$rng = $d.Content\n$found = $rng.Find.Execute(\"null}\")\n$xml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r w:rsidR=\"00D500A1\"><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>null</w:t></w:r><w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$rng.InsertOoxml($xml)\nWrite-Output $d.Content.Text\n"}
